$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alumnos")

# Row 28: topic changed from "Comparación de enfoques..." to "Retrospectiva"
$ws.Range("D28").Value = "Retrospectiva"

# Row 28: instructions text extended with a note about it being used as a makeup
$ws.Range("E28").Value = "Deben traer para trabajar en clase un cuadro comparativo de los 3 enfoques (ver TP 14 en la guía de prácticos). Es evaluable. Se usará como recuperatorio de los TP conceptuales, para los grupos que lo necesiten"

# Row 28: new recorded-class link in the Material column
$ws.Hyperlinks.Add($ws.Range("F28"), "https://youtu.be/u1bmaI4bEaU") | Out-Null
$ws.Range("F28").Value = "Clase Grabada Retrospectiva"
$ws.Range("F28").Font.Name = "Calibri"
$ws.Range("F28").Font.Color = 13391121
$ws.Range("F28").Font.Underline = 2

# Row 28 needs to grow to fit the wrapped instructions text
$ws.Rows.Item(28).RowHeight = 54

# Row 29: new recorded-class link in the "Acceso a clase grabada" column
$ws.Hyperlinks.Add($ws.Range("G29"), "https://youtu.be/u1bmaI4bEaU") | Out-Null
$ws.Range("G29").Value = "Clase Grabada Practico 13"
$ws.Range("G29").Font.Name = "Calibri"
$ws.Range("G29").Font.Color = 13391121
$ws.Range("G29").Font.Underline = 2
